$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (existing row): D2 changes from text "07678447952" to numeric 7678447952 ---
$ws.Cells.Item(2, 4).Value = 7678447952

# --- Row 3 (new) ---
$ws.Cells.Item(3, 1).Value = "AAYUSH MISHRA"
$ws.Cells.Item(3, 2).Value = "10323210240@stu.smuniversity.ac.in"
$ws.Cells.Item(3, 3).Value = "SRM"
$ws.Cells.Item(3, 4).Value = 7678447952
$ws.Cells.Item(3, 5).Value = "Visitor Management System"
$ws.Cells.Item(3, 6).Value = ""

# --- Row 4 (new) ---
$ws.Cells.Item(4, 1).Value = "Aayush Mishra"
$ws.Cells.Item(4, 2).Value = "aayushmishra82017@gmail.com"
$ws.Cells.Item(4, 3).Value = "NGAT Innovations "
$ws.Cells.Item(4, 4).Value = 7678447952
$ws.Cells.Item(4, 5).Value = "Bill Organiser, Visitor Management System, Stock Manager and Tracker"
$ws.Cells.Item(4, 6).Value = ""

# --- Row 5 (new) ---
$ws.Cells.Item(5, 1).Value = "Raghav "
$ws.Cells.Item(5, 2).Value = "aayushmishra82018@gmail.com"
$ws.Cells.Item(5, 3).Value = "Ashoka"
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "12345676i7"
$ws.Cells.Item(5, 5).Value = "Stock Manager and Tracker, Machine Fault and Maintenance Detector, Automated Product Fault Detector"
$ws.Cells.Item(5, 6).Value = "i need a personalised server space "

# --- Row 6 (new) ---
$ws.Cells.Item(6, 1).Value = "Raghav Saini "
$ws.Cells.Item(6, 2).Value = "raghavsaini9560@gmail.com"
$ws.Cells.Item(6, 3).Value = "NGAT"
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "7701970160"
$ws.Cells.Item(6, 5).Value = "Data Handler, Production Line Record Keeper, Automated Product Fault Detector, Machine Line Organiser"
$ws.Cells.Item(6, 6).Value = "I need my Server space"

# --- Row 7 (new) ---
$ws.Cells.Item(7, 1).Value = "Aayush Mishra"
$ws.Cells.Item(7, 2).Value = "aayushmishra82018@gmail.com"
$ws.Cells.Item(7, 3).Value = "SRM"
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "07678447952"
$ws.Cells.Item(7, 5).Value = "Stock Manager and Tracker"
$ws.Cells.Item(7, 6).Value = ""
